$wb = $excel.ActiveWorkbook

# Update Sheet2: A4 changes from 5 to 4 (B4 formula A4*7 recalculates to 28)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A4").Value = 4

# Make Sheet2 the active sheet/tab
$ws2.Activate()
